# Bank_Data.xlsx update ("Updated all Bank files")
#
# - A2 (Branch_Num) was re-keyed from the old placeholder date-like value
#   20240519 to the real branch number 1005151.
# - The cell cursor / selection that was left on D3 is moved back to A3
#   (i.e. just below the data that was actually edited in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Branch_Num value in A2.
$ws.Range("A2").Value = 1005151

# Move the active selection to A3 (was D3).
$ws.Range("A3").Select()
